$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.899.83"
$ws.Range("E2").Value = "  -7.57%  "

$ws.Range("D3").Value = "3.726.90"
$ws.Range("E3").Value = "  -6.31%  "

$ws.Range("E4").Value = "  +0.23%  "

$ws.Range("D5").Value = "574.55"
$ws.Range("E5").Value = "  -5.93%  "

$ws.Range("D6").Value = "174.32"
$ws.Range("E6").Value = "  +3.18%  "

$ws.Range("D7").Value = "3.717.39"
$ws.Range("E7").Value = "  -6.30%  "

$ws.Range("D8").Value = "0.630"
$ws.Range("E8").Value = "  -7.42%  "

$ws.Range("D9").Value = "0.998"
$ws.Range("E9").Value = "  -0.16%  "

$ws.Range("D10").Value = "0.711"
$ws.Range("E10").Value = "  -8.91%  "

$ws.Range("D11").Value = "0.165"
$ws.Range("E11").Value = "  -10.56%  "

$ws.Range("D12").Value = "53.02"
$ws.Range("E12").Value = "  -5.62%  "

$ws.Range("D13").Value = "0.0000298"
$ws.Range("E13").Value = "  -10.93%  "

$ws.Range("D14").Value = "10.63"
$ws.Range("E14").Value = "  -5.68%  "

$ws.Range("D15").Value = "4.342.12"
$ws.Range("E15").Value = "  -5.76%  "

$ws.Range("D16").Value = "3.747.94"
$ws.Range("E16").Value = "  -5.45%  "

$ws.Range("D17").Value = "19.38"
$ws.Range("E17").Value = "  -5.57%  "

$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").Value = "0.127"
$ws.Range("E18").Value = "  -2.89%  "

$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").Value = "12.97"
$ws.Range("E19").Value = "  -8.71%  "

$ws.Range("B20").Value = "Polygon"
$ws.Range("C20").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D20").Value = "1.14"
$ws.Range("E20").Value = "  -8.01%  "

$ws.Range("D21").Value = "68.035.80"
$ws.Range("E21").Value = "  -7.14%  "

$ws.Range("D22").Value = "408.19"
$ws.Range("E22").Value = "  -12.50%  "

$ws.Range("D23").Value = "4.51"
$ws.Range("E23").Value = "  -6.65%  "

$ws.Range("D24").Value = "88.57"
$ws.Range("E24").Value = "  -8.58%  "

$ws.Range("D25").Value = "3.07"
$ws.Range("E25").Value = "  -10.06%  "

$ws.Range("D26").Value = "12.81"
$ws.Range("E26").Value = "  -9.73%  "

$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D27").Value = "10.66"
$ws.Range("E27").Value = "  -3.26%  "

$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").Value = "3.85"
$ws.Range("E28").Value = "  -8.59%  "

$ws.Range("E29").Value = "  +0.32%  "

$ws.Range("D30").Value = "9.56"
$ws.Range("E30").Value = "  -9.19%  "

$ws.Range("D31").Value = "32.93"
$ws.Range("E31").Value = "  -9.28%  "

$ws.Range("D32").Value = "7.74"
$ws.Range("E32").Value = "  -1.72%  "

$ws.Range("D33").Value = "12.65"
$ws.Range("E33").Value = "  -9.08%  "

$ws.Range("D34").Value = "0.117"
$ws.Range("E34").Value = "  -9.18%  "

$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").Value = "65.32"
$ws.Range("E35").Value = "  -7.57%  "

$ws.Range("B36").Value = "InjectiveProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D36").Value = "43.56"
$ws.Range("E36").Value = "  -9.02%  "

$ws.Range("D37").Value = "604.84"
$ws.Range("E37").Value = "  -6.69%  "

$ws.Range("D38").Value = "0.0₃0908"
$ws.Range("E38").Value = "  -13.09%  "

$ws.Range("D39").Value = "0.400"
$ws.Range("E39").Value = "  -6.84%  "

$ws.Range("E40").Value = "  +0.22%  "

$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  +0.37%  "

$ws.Range("D42").Value = "0.136"
$ws.Range("E42").Value = "  -6.56%  "

$ws.Range("B43").Value = "ThetaToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D43").Value = "3.05"
$ws.Range("E43").Value = "  -10.00%  "

$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").Value = "3.00"
$ws.Range("E44").Value = "  -5.15%  "

$ws.Range("D45").Value = "0.0441"
$ws.Range("E45").Value = "  -8.70%  "

$ws.Range("D46").Value = "2.58"
$ws.Range("E46").Value = "  +0.75%  "

$ws.Range("D47").Value = "9.33"
$ws.Range("E47").Value = "  -10.78%  "

$ws.Range("B48").Value = "WEMIXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D48").Value = "2.73"
$ws.Range("E48").Value = "  -12.96%  "

$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").Value = "0.135"
$ws.Range("E49").Value = "  -9.68%  "

$ws.Range("D50").Value = "2.754.46"
$ws.Range("E50").Value = "  -2.54%  "

$ws.Range("E51").Value = "  -8.02%  "
